# The commit simply renames the two logo pictures that live in the
# document's headers/footers:
#   - the BTec logo inline picture  (descr "BTec_Logo-Orange"):      image1.jpg -> image2.jpg
#   - the Pearson logo inline picture (descr "...PearsonLogo.png"):  image2.png -> image1.png
#
# NB: InlineShape.Name's getter only ever reflects a value set earlier in
# this same session (it does not read back the pre-existing wp:docPr/@name
# from the file), so we must not branch on the current .Name. Instead we
# identify each picture by its (file-accurate) AlternativeText/descr and
# then assign the exact target name called for by the commit.

$d = $word.ActiveDocument

function Set-LogoShapeName($shape) {
    if ($null -eq $shape) { return }
    $descr = $shape.AlternativeText

    if ($descr -eq "BTec_Logo-Orange") {
        $shape.Name = "image2.jpg"
    } elseif ($descr -like "*PearsonLogo.png") {
        $shape.Name = "image1.png"
    }
}

foreach ($sec in $d.Sections) {
    # wdHeaderFooterPrimary=1, wdHeaderFooterFirstPage=2, wdHeaderFooterEvenPages=3
    for ($hi = 1; $hi -le 3; $hi++) {
        $hf = $sec.Headers.Item($hi)
        $cnt = $hf.Range.InlineShapes.Count
        for ($i = 1; $i -le $cnt; $i++) {
            Set-LogoShapeName $hf.Range.InlineShapes.Item($i)
        }
    }
    for ($fi = 1; $fi -le 3; $fi++) {
        $hf = $sec.Footers.Item($fi)
        $cnt = $hf.Range.InlineShapes.Count
        for ($i = 1; $i -le $cnt; $i++) {
            Set-LogoShapeName $hf.Range.InlineShapes.Item($i)
        }
    }
}
